$d = $word.ActiveDocument

$map = @(
    @("67×52=3484", "11×34=374"),
    @("25×97=2425", "58×90=5220"),
    @("84×93=7812", "91×78=7098"),
    @("11×49=539",  "70×50=3500"),
    @("48×79=3792", "27×82=2214"),
    @("81×47=3807", "87×49=4263"),
    @("46×52=2392", "40×57=2280"),
    @("37×18=666",  "63×37=2331"),
    @("43×11=473",  "52×27=1404"),
    @("13×19=247",  "74×35=2590"),
    @("91×40=3640", "80×54=4320"),
    @("80×74=5920", "76×89=6764"),
    @("36×34=1224", "71×55=3905"),
    @("23×64=1472", "99×92=9108"),
    @("96×63=6048", "18×13=234"),
    @("97×60=5820", "93×44=4092"),
    @("33×84=2772", "73×42=3066"),
    @("71×81=5751", "55×91=5005"),
    @("29×94=2726", "79×86=6794"),
    @("17×82=1394", "95×19=1805"),
    @("64×20=1280", "67×80=5360"),
    @("37×45=1665", "62×76=4712"),
    @("85×69=5865", "71×36=2556"),
    @("99×85=8415", "58×72=4176"),
    @("46×11=506",  "75×90=6750")
)

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
